$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Targeted single-cell updates in existing rows ---
$ws.Cells.Item(55, 17).Value = 0      # Q55: 2 -> 0
$ws.Cells.Item(450, 15).Value = 2     # O450: 0 -> 2
$ws.Cells.Item(452, 15).Value = 1     # O452: 0 -> 1
$ws.Cells.Item(452, 18).Value = 0     # R452: blank -> 0
$ws.Cells.Item(453, 18).Value = 0     # R453: blank -> 0

# --- New weekly rows 454-465 ---
$dateNumFmt = $ws.Cells.Item(453, 1).NumberFormat()

$newRows = @(
    @(45474, 82.16000366210938, 82.19999694824219, 77.84999847412109, 81.19000244140625, 81.19000244140625, 197627427, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 81.18000030517578, 81.34999847412109, 77.61000061035156, 78.26999664306641, 78.26999664306641, 122112384, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 78.30000305175781, 78.80000305175781, 75.84999847412109, 76.01999664306641, 76.01999664306641, 74564791, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 75.94000244140625, 78.04000091552734, 72.40000152587891, 74.48000335693359, 74.48000335693359, 186272494, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 73.69000244140625, 76.45999908447266, 73.05000305175781, 74.30999755859375, 74.30999755859375, 199811094, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 2),
    @(45509, 72.98999786376953, 73.69999694824219, 71.55999755859375, 72.86000061035156, 72.86000061035156, 137392398, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 72.51000213623047, 72.77999877929688, 70.43000030517578, 71.95999908447266, 71.95999908447266, 119738948, 2024, 8, 12, 0, 0, 0, 33, 2, 0, 0),
    @(45523, 72, 75.79000091552734, 71.59999847412109, 74.41999816894531, 74.41999816894531, 139457021, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 74.87000274658203, 75.05000305175781, 72.94000244140625, 73.83999633789062, 73.83999633789062, 106216830, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 73.65000152587891, 75.69999694824219, 73, 73.66000366210938, 73.66000366210938, 125773472, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 73.58999633789062, 74.18000030517578, 71.36000061035156, 73.41999816894531, 73.41999816894531, 147102165, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(45551, 73.63999938964844, 74.62999725341797, 71.16000366210938, 72.83000183105469, 72.83000183105469, 175077866, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0)
)

$startRow = 454
foreach ($row in $newRows) {
    $ws.Cells.Item($startRow, 1).Value = $row[0]
    $ws.Cells.Item($startRow, 1).NumberFormat = $dateNumFmt
    for ($i = 1; $i -lt $row.Count; $i++) {
        $ws.Cells.Item($startRow, $i + 1).Value = $row[$i]
    }
    $startRow = $startRow + 1
}

